$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget Worksheet")

# --- Data edits, applied in the same order the author made them (so any
#     newly-created shared strings line up with the expected indices) ---

# C31: 156 -> 160
$ws.Range("C31").Value = 160

# C36: 496 -> 500
$ws.Range("C36").Value = 500

# D31 gains a right-aligned, indented number format as part of this edit
$ws.Range("D31").HorizontalAlignment = -4152
$ws.Range("D31").IndentLevel = 1

# E25: "US Aquaculture Society" -> "SAFS FINS"
$ws.Range("E25").Value = "SAFS FINS"

# D31: 400 -> 375
$ws.Range("D31").Value = 375

# D25: 100 -> 121
$ws.Range("D25").Value = 121

# F11 / F25 / F31: new "Pending" funding-approved notes
$ws.Range("F11").Value = "Pending"
$ws.Range("F25").Value = "Pending"
$ws.Range("F31").Value = "Pending"

# G11 / G16 / G13 / G25 / G31: new explanation-of-expense notes
$ws.Range("G11").Value = "Student registration fee"
$ws.Range("G16").Value = "Transit within Las Vegas"
$ws.Range("G13").Value = "Airfare SEA <-> Las Vegas"
$ws.Range("G25").Value = "Meals not provided by AA"

# G20 total is re-entered (the shared SUM formula over G13:G19 is unchanged in effect)
$ws.Range("G20").Formula = "=SUM(G13:G19)"

$ws.Range("G31").Value = "Hotel during conference"

# --- Selection moves to F36 (last place the author clicked before saving) ---
$ws.Range("F36").Select()

# --- Page margins reset to Excel's normal defaults on both sheets ---
foreach ($sheetName in @("Budget Worksheet", "Sample Worksheet")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}
